$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the two cells that were missing at the end of row 10 ---
$ws.Cells.Item(10, 24).Value = -0.16000300000000323
$ws.Cells.Item(10, 25).Value = "Down"

# --- Append the new row 11 of data ---
$ws.Cells.Item(11, 1).Value = 42654.894479166665
$ws.Cells.Item(11, 2).Value = 9
$ws.Cells.Item(11, 3).Value = "Buy"
$ws.Cells.Item(11, 4).Value = 6
$ws.Cells.Item(11, 5).Value = 14166
$ws.Cells.Item(11, 6).Value = 2569
$ws.Cells.Item(11, 7).Value = 52
$ws.Cells.Item(11, 8).Value = 43
$ws.Cells.Item(11, 9).Value = 74
$ws.Cells.Item(11, 10).Value = 25
$ws.Cells.Item(11, 11).Value = 17319
$ws.Cells.Item(11, 12).Value = 264
$ws.Cells.Item(11, 13).Value = 215
$ws.Cells.Item(11, 14).Value = 73
$ws.Cells.Item(11, 15).Value = 25
$ws.Cells.Item(11, 16).Value = "Noun"
$ws.Cells.Item(11, 17).Value = 28.689659976213832
$ws.Cells.Item(11, 18).Value = 0.84
$ws.Cells.Item(11, 19).Value = -0.012500000000000001
$ws.Cells.Item(11, 20).Value = -0.026100000000000002
$ws.Cells.Item(11, 21).Value = 14.56
$ws.Cells.Item(11, 22).Value = "N/A"
$ws.Cells.Item(11, 23).Value = 0

# Columns S (19) and T (20) use a percentage number format on existing
# rows (style index 2) that isn't carried automatically onto a brand new
# row, so copy it explicitly from the row above.
$ws.Cells.Item(11, 19).NumberFormat = $ws.Cells.Item(10, 19).NumberFormat
$ws.Cells.Item(11, 20).NumberFormat = $ws.Cells.Item(10, 20).NumberFormat
